{"js": "// Update the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n// in the single table: each cell holds 5 lines (separated by manual line\n// breaks) -- \"A x B\", the spaced-out digits of B, a \"----\" divider, and\n// the two digits of A each followed by \"|    |\".\n//\n// final[row][col] gives the five lines for that cell in their new\n// (post-edit) state -- including lines that happen to be unchanged, so\n// each touched cell is rewritten consistently in one shot.\nconst finalValues = [\n  [\n    [\"57 x 39\", \"  3    9\", \"  ----\", \"5|    |\", \"7|    |\"],\n    [\"66 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"6|    |\"],\n    [\"92 x 30\", \"  3    0\", \"  ----\", \"9|    |\", \"2|    |\"],\n  ],\n  [\n    [\"22 x 20\", \"  2    0\", \"  ----\", \"2|    |\", \"2|    |\"],\n    [\"40 x 79\", \"  7    9\", \"  ----\", \"4|    |\", \"0|    |\"],\n    [\"58 x 39\", \"  3    9\", \"  ----\", \"5|    |\", \"8|    |\"],\n  ],\n  [\n    [\"49 x 34\", \"  3    4\", \"  ----\", \"4|    |\", \"9|    |\"],\n    [\"98 x 83\", \"  8    3\", \"  ----\", \"9|    |\", \"8|    |\"],\n    [\"95 x 84\", \"  8    4\", \"  ----\", \"9|    |\", \"5|    |\"],\n  ],\n  [\n    [\"28 x 26\", \"  2    6\", \"  ----\", \"2|    |\", \"8|    |\"],\n    [\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"],\n    [\"75 x 50\", \"  5    0\", \"  ----\", \"7|    |\", \"5|    |\"],\n  ],\n  [\n    [\"28 x 72\", \"  7    2\", \"  ----\", \"2|    |\", \"8|    |\"],\n    [\"28 x 18\", \"  1    8\", \"  ----\", \"2|    |\", \"8|    |\"],\n    [\"26 x 64\", \"  6    4\", \"  ----\", \"2|    |\", \"6|    |\"],\n  ],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let row = 0; row < finalValues.length; row++) {\n  for (let col = 0; col < finalValues[row].length; col++) {\n    const cell = table.getCellOrNullObject(row, col);\n    cell.load(\"isNullObject\");\n    await context.sync();\n    if (cell.isNullObject) {\n      continue;\n    }\n    const lines = finalValues[row][col];\n    const cellRange = cell.body.getRange();\n    cellRange.insertText(lines.join(\"\\v\"), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n# in the single table: each cell holds 5 lines (separated by manual line\n# breaks, Chr(11)) -- \"A x B\", the spaced-out digits of B, a \"----\"\n# divider, and the two digits of A each followed by \"|    |\".\n#\n# $final holds, for every (row, col), the five lines in their new\n# (post-edit) state -- including lines that happen to be unchanged, so\n# each cell is rewritten consistently in one shot.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$nl = [char]11\n\n$final = @(\n    @(\n        @(\"57 x 39\", \"  3    9\", \"  ----\", \"5|    |\", \"7|    |\"),\n        @(\"66 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"6|    |\"),\n        @(\"92 x 30\", \"  3    0\", \"  ----\", \"9|    |\", \"2|    |\")\n    ),\n    @(\n        @(\"22 x 20\", \"  2    0\", \"  ----\", \"2|    |\", \"2|    |\"),\n        @(\"40 x 79\", \"  7    9\", \"  ----\", \"4|    |\", \"0|    |\"),\n        @(\"58 x 39\", \"  3    9\", \"  ----\", \"5|    |\", \"8|    |\")\n    ),\n    @(\n        @(\"49 x 34\", \"  3    4\", \"  ----\", \"4|    |\", \"9|    |\"),\n        @(\"98 x 83\", \"  8    3\", \"  ----\", \"9|    |\", \"8|    |\"),\n        @(\"95 x 84\", \"  8    4\", \"  ----\", \"9|    |\", \"5|    |\")\n    ),\n    @(\n        @(\"28 x 26\", \"  2    6\", \"  ----\", \"2|    |\", \"8|    |\"),\n        @(\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"),\n        @(\"75 x 50\", \"  5    0\", \"  ----\", \"7|    |\", \"5|    |\")\n    ),\n    @(\n        @(\"28 x 72\", \"  7    2\", \"  ----\", \"2|    |\", \"8|    |\"),\n        @(\"28 x 18\", \"  1    8\", \"  ----\", \"2|    |\", \"8|    |\"),\n        @(\"26 x 64\", \"  6    4\", \"  ----\", \"2|    |\", \"6|    |\")\n    )\n)\n\nfor ($row = 0; $row -lt $final.Count; $row++) {\n    for ($col = 0; $col -lt $final[$row].Count; $col++) {\n        $lines = $final[$row][$col]\n        $cell = $t.Cell($row + 1, $col + 1)\n        $cell.Range.Text = [string]::Join($nl, $lines)\n    }\n}\n"}
